$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-26 12:46:08"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-08-26 12:45:58"
$zhcn.Range("K4").Value = "2016-08-26 12:46:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K4").Value = "2016-08-26 12:46:37"
